# Regenerate the "K" column (column G) values on Sheet1.
# The save_data sheet used to store a "Strike#" derived value in column G
# ("K"); this re-run recalculates/rewrites those values (s_vals) for each
# row while leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K (column G) value, as produced by the regenerated
# save_data computation.
$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 2
    14 = 1
    15 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
